$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 6: Median of Array -----
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "Binary S 14"
$ws.Range("D6").Value = "Median of Array"

$e6 = $ws.Range("E6")
$ws.Hyperlinks.Add($e6, "https://www.scaler.com/academy/mentee-dashboard/class/30364/homework/problems/198/?navref=cl_pb_nv_tb") | Out-Null
$e6.Style = "Hyperlink"
$e6.Style.VerticalAlignment = -4108
$e6.Style.WrapText = $true

$f6 = $ws.Range("F6")
$ws.Hyperlinks.Add($f6, "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/4_median_of_array.java", "", "", "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/4_median_of_array.java") | Out-Null
$f6.Value = "dsa/4_median_of_array.java at main · ankurnecessary/dsa · GitHub"
$f6.Style = "Hyperlink"
$f6.Style.WrapText = $true

$ws.Rows.Item(6).RowHeight = 72

# ----- Row 7: Just greater number -----
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Binary S"
$ws.Range("D7").Value = "Just greater number"
$ws.Range("E7").Value = "No Link"

$f7 = $ws.Range("F7")
$ws.Hyperlinks.Add($f7, "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/5_just_greater_number.java", "", "", "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/5_just_greater_number.java") | Out-Null
$f7.Value = "dsa/5_just_greater_number.java at main · ankurnecessary/dsa · GitHub"
$f7.Style = "Hyperlink"
$f7.Style.WrapText = $true

$ws.Rows.Item(7).RowHeight = 43.2

# ----- Selection -----
$ws.Range("F7").Select() | Out-Null

Write-Host "done"
